$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "Was there a shift in market sentiment when Bitcoin was created? ..."
#    Split the sentence so a new "(popularity) " run is inserted right after
#    "market sentiment ".
# ---------------------------------------------------------------------------
$p10 = $d.Paragraphs(10)
$r10 = $p10.Range
$found = $r10.Find.Execute("market sentiment ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r10.Collapse(0)
$r10.InsertAfter("(popularity) ")

# ---------------------------------------------------------------------------
# 2. "What are benefits to using crypto compared to standard currency?" ->
#    "Which do people prefer more? Class vote "
#    "How will Bitcoin be regulated?" ->
#    "Which is more volatile/which one is a riskier investment?"
#    Plus a brand new bullet: "Which countries own the most of each currency?"
# ---------------------------------------------------------------------------
$p12 = $d.Paragraphs(12)
$p12.Range.Text = "Which do people prefer more? Class vote "

$p13 = $d.Paragraphs(13)
$p13.Range.Text = "Which is more volatile/which one is a riskier investment?"

$newPar = $d.Paragraphs(13).Range.InsertParagraphAfter()
$d.Paragraphs(14).Range.Text = "Which countries own the most of each currency?"

# ---------------------------------------------------------------------------
# 3. "TBD" -> "Data - TJ and Key" followed by five new task bullets.
# ---------------------------------------------------------------------------
$tbdIndex = $d.Paragraphs.Count
$pTbd = $d.Paragraphs($tbdIndex)
$pTbd.Range.Text = "Data - TJ and Key"

$d.Paragraphs($tbdIndex).Range.InsertParagraphAfter()
$tbdIndex = $tbdIndex + 1
$d.Paragraphs($tbdIndex).Range.Text = "SQL – Key"

$d.Paragraphs($tbdIndex).Range.InsertParagraphAfter()
$tbdIndex = $tbdIndex + 1
$d.Paragraphs($tbdIndex).Range.Text = "Flask API - TJ"

$d.Paragraphs($tbdIndex).Range.InsertParagraphAfter()
$tbdIndex = $tbdIndex + 1
$d.Paragraphs($tbdIndex).Range.Text = "Bar charts/HTML – Wipawadee "

$d.Paragraphs($tbdIndex).Range.InsertParagraphAfter()
$tbdIndex = $tbdIndex + 1
$d.Paragraphs($tbdIndex).Range.Text = "Line charts/CSS - Shloka"

$d.Paragraphs($tbdIndex).Range.InsertParagraphAfter()
$tbdIndex = $tbdIndex + 1
$d.Paragraphs($tbdIndex).Range.Text = "Heatmap - Everybody"
